$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 62 (shifts existing rows 62:86 down to 63:87,
# pushing the sheet's used range from R86 to R87), copying formatting from
# the row above (keeps the date-format style on column D).
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new weekly price record.
$ws.Range("A62").Value = 4
$ws.Range("B62").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C62").Value = 'Los Lagos'
$ws.Range("D62").Value = 45027
$ws.Range("E62").Value = 10
$ws.Range("F62").Value = 100112043
$ws.Range("G62").Value = 'Pepino dulce'
$ws.Range("H62").Value = 'Cultivar IV Región'
$ws.Range("I62").Value = 'Primera'
$ws.Range("J62").Value = 80
$ws.Range("K62").Value = 18000
$ws.Range("L62").Value = 20000
$ws.Range("M62").Value = 19000
$ws.Range("N62").Value = '$/bandeja 18 kilos'
$ws.Range("O62").Value = 'Provincia de Limarí'
$ws.Range("P62").Value = 1056
$ws.Range("Q62").Value = 18
$ws.Range("R62").Value = 'Hortaliza'
